$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Exported On:" timestamp shown near the top of the report.
$ws.Range("A2").Value = "Exported On: 06/23/2025 05:48 PM"

# Update the per-row Total Cost (F) and Quantity (K) figures that changed
# between the two exports.
$rowUpdates = @(
    @{ Row = 12;  F = 0.0;     K = -2.0 },
    @{ Row = 13;  F = 0.0;     K = 0.0 },
    @{ Row = 16;  F = 13.42;   K = 2.0 },
    @{ Row = 20;  F = 1294.7;  K = 2.0 },
    @{ Row = 21;  F = 875.26;  K = 2.0 },
    @{ Row = 36;  F = 51.06;   K = 46.0 },
    @{ Row = 38;  F = 99.88;   K = 22.0 },
    @{ Row = 41;  F = 18.36;   K = 9.0 },
    @{ Row = 54;  F = 14.94;   K = 6.0 },
    @{ Row = 59;  F = 16.24;   K = 7.0 },
    @{ Row = 92;  F = 0.0;     K = 0.0 },
    @{ Row = 101; F = 0.0;     K = 0.0 },
    @{ Row = 102; F = 40.32;   K = 6.0 },
    @{ Row = 130; F = 80.0;    K = 4.0 },
    @{ Row = 131; F = 85.52;   K = 8.0 },
    @{ Row = 139; F = 42.72;   K = 4.0 },
    @{ Row = 144; F = 38.46;   K = 3.0 },
    @{ Row = 147; F = 0.0;     K = 0.0 },
    @{ Row = 157; F = 17.1;    K = 2.0 },
    @{ Row = 176; F = 14.08;   K = 11.0 },
    @{ Row = 183; F = 35.28;   K = 3.0 },
    @{ Row = 185; F = 0.0;     K = 0.0 }
)

foreach ($u in $rowUpdates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 11).Value = $u.K
}

# Update the two "Total" summary rows at the bottom of the report.
$ws.Cells.Item(193, 6).Value = 22066.06
$ws.Cells.Item(193, 11).Value = 731.0
$ws.Cells.Item(195, 6).Value = 22066.06
$ws.Cells.Item(195, 11).Value = 731.0
